$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

# --- Refresh the panel_query_time-derived "time_taken" column on the "data" sheet ---
$data.Range("F2").Value = "2021-10-05 14:35:52.487075"
$data.Range("F3").Value = "2021-10-05 14:35:52.487083"
$data.Range("F4").Value = "2021-10-05 14:35:52.487086"
$data.Range("F5").Value = "2021-10-05 14:35:52.487089"
$data.Range("F6").Value = "2021-10-05 14:35:52.487092"
$data.Range("F7").Value = "2021-10-05 14:35:52.487095"
$data.Range("F8").Value = "2021-10-05 14:35:52.487097"
$data.Range("F9").Value = "2021-10-05 14:35:52.487100"

# --- Add a new "metadata" worksheet right after "data" ---
$metadata = $wb.Worksheets.Add($null, $data)
$metadata.Name = "metadata"

# Header row
$metadata.Range("B1").Value = "data_name"
$metadata.Range("C1").Value = "data_id"
$metadata.Range("D1").Value = "data_version"
$metadata.Range("E1").Value = "data_version_created"
$metadata.Range("F1").Value = "panel_query_time"
$metadata.Range("G1").Value = "panel_get_request"

# Match the bold/bordered header style used on the "data" sheet
$data.Range("B1:F1").Copy()
$metadata.Range("B1:G1").PasteSpecial(-4122)  # xlPasteFormats

# Data row
$metadata.Range("A2").Value = 0
$data.Range("A2").Copy()
$metadata.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

$metadata.Range("B2").Value = "Tubulinopathies"
$metadata.Range("C2").Value = 21

# data_version must stay the literal text "1.0" rather than become the number 1
$metadata.Range("D2").Formula = '="1.0"'
$metadata.Range("D2").Copy()
$metadata.Range("D2").PasteSpecial(-4163)  # xlPasteValues

$metadata.Range("E2").Value = "2021-01-16T10:10:49.613318Z"
$metadata.Range("F2").Value = "2021-10-05 14:35:52.483238"
$metadata.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/21/?format=json"

$excel.CutCopyMode = 0

# Keep "data" as the active sheet/tab, as in the original workbook
$data.Activate()
